$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells are formatted as Text so the literal numeric/percent-looking
# strings are preserved exactly as authored (matching the source inline-string cells)
# rather than being auto-coerced into Number/Percentage values by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.10%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.229"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.98%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07669"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.70%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.643"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.08%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9162"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.84%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.432"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1250"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "15.07%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1828"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.77%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09164"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.40%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04259"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.28%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001263"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.67%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005755"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.62%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.308"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.28%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.321"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "11.47%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1383"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.38%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2717"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.35%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04069"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.48%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.32%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004288"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.72%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001273"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.11%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02477"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.56%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05296"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.36%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007847"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.12%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.20%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006883"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.47%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001915"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.86%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007644"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.73%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3057"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.61%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006733"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.92%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1699"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "424.57%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-26.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
